$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (preserving exact string formatting,
# e.g. "42.752.19" or "116.97") without Excel auto-converting it to a
# number. Temporarily flips the cell to Text format, assigns the value,
# then restores the "Normal" cell style so no lasting formatting change
# is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.748.86"
Set-TextValue $ws.Range("E2") "  -2.31%  "

Set-TextValue $ws.Range("D3") "2.248.19"
Set-TextValue $ws.Range("E3") "  -1.82%  "

Set-TextValue $ws.Range("E4") "  +0.04%  "

Set-TextValue $ws.Range("D5") "116.53"
Set-TextValue $ws.Range("E5") "  +0.91%  "

Set-TextValue $ws.Range("D6") "297.33"
Set-TextValue $ws.Range("E6") "  +11.65%  "

Set-TextValue $ws.Range("D7") "0.633"
Set-TextValue $ws.Range("E7") "  -1.67%  "

Set-TextValue $ws.Range("E8") "  -0.07%  "

Set-TextValue $ws.Range("D9") "0.621"
Set-TextValue $ws.Range("E9") "  +1.00%  "

Set-TextValue $ws.Range("D10") "46.53"
Set-TextValue $ws.Range("E10") "  -2.09%  "

Set-TextValue $ws.Range("E11") "  -0.34%  "

Set-TextValue $ws.Range("D12") "9.15"
Set-TextValue $ws.Range("E12") "  -0.45%  "

Set-TextValue $ws.Range("E13") "  -2.75%  "

Set-TextValue $ws.Range("D14") "15.48"
Set-TextValue $ws.Range("E14") "  +0.35%  "

Set-TextValue $ws.Range("D15") "0.894"
Set-TextValue $ws.Range("E15") "  +2.08%  "

Set-TextValue $ws.Range("D16") "2.588.12"
Set-TextValue $ws.Range("E16") "  -1.85%  "

Set-TextValue $ws.Range("D17") "2.261.92"
Set-TextValue $ws.Range("E17") "  -1.46%  "

Set-TextValue $ws.Range("D18") "42.846.22"
Set-TextValue $ws.Range("E18") "  -1.96%  "

Set-TextValue $ws.Range("D19") "7.56"
Set-TextValue $ws.Range("E19") "  +11.49%  "

Set-TextValue $ws.Range("E20") "  -1.87%  "

Set-TextValue $ws.Range("D21") "73.95"
Set-TextValue $ws.Range("E21") "  +1.97%  "

Set-TextValue $ws.Range("D22") "3.47"
Set-TextValue $ws.Range("E22") "  +19.97%  "

Set-TextValue $ws.Range("D23") "2.37"
Set-TextValue $ws.Range("E23") "  -2.90%  "

Set-TextValue $ws.Range("D24") "233.38"
Set-TextValue $ws.Range("E24") "  -1.45%  "

Set-TextValue $ws.Range("D25") "9.43"
Set-TextValue $ws.Range("E25") "  -0.67%  "

Set-TextValue $ws.Range("D26") "12.30"
Set-TextValue $ws.Range("E26") "  +5.97%  "

Set-TextValue $ws.Range("D27") "1.00"
Set-TextValue $ws.Range("E27") "  -1.81%  "

Set-TextValue $ws.Range("D28") "40.27"
Set-TextValue $ws.Range("E28") "  -3.56%  "

Set-TextValue $ws.Range("D31") "175.92"
Set-TextValue $ws.Range("E31") "  +1.06%  "

Set-TextValue $ws.Range("D32") "21.34"
Set-TextValue $ws.Range("E32") "  -1.87%  "

Set-TextValue $ws.Range("D33") "0.0912"
Set-TextValue $ws.Range("E33") "  +0.19%  "

Set-TextValue $ws.Range("D34") "4.61"
Set-TextValue $ws.Range("E34") "  +16.08%  "

Set-TextValue $ws.Range("E35") "  -0.53%  "

Set-TextValue $ws.Range("E36") "  -1.74%  "

Set-TextValue $ws.Range("D37") "4.77"
Set-TextValue $ws.Range("E37") "  +1.26%  "

Set-TextValue $ws.Range("D38") "0.0376"
Set-TextValue $ws.Range("E38") "  -1.37%  "

Set-TextValue $ws.Range("E39") "  +0.71%  "

Set-TextValue $ws.Range("D40") "2.62"
Set-TextValue $ws.Range("E40") "  +2.41%  "

Set-TextValue $ws.Range("D41") "72.84"
Set-TextValue $ws.Range("E41") "  -2.33%  "

Set-TextValue $ws.Range("D42") "0.239"
Set-TextValue $ws.Range("E42") "  +1.05%  "

Set-TextValue $ws.Range("D43") "13.56"
Set-TextValue $ws.Range("E43") "  -5.21%  "

Set-TextValue $ws.Range("E44") "  +0.15%  "

Set-TextValue $ws.Range("E45") "  -1.82%  "

Set-TextValue $ws.Range("E46") "  -7.36%  "

Set-TextValue $ws.Range("D47") "1.33"
Set-TextValue $ws.Range("E47") "  +2.95%  "

Set-TextValue $ws.Range("D48") "106.83"
Set-TextValue $ws.Range("E48") "  +6.45%  "

Set-TextValue $ws.Range("D49") "8.64"
Set-TextValue $ws.Range("E49") "  +0.40%  "

Set-TextValue $ws.Range("E50") "  -0.85%  "

Set-TextValue $ws.Range("E51") "  +6.00%  "

# Rows 29 & 30: Toncoin and WEMIXToken swapped positions, each with updated values
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "2.23"
Set-TextValue $ws.Range("E29") "  -0.89%  "

$ws.Range("B30").Value = "WEMIXToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D30") "3.28"
Set-TextValue $ws.Range("E30") "  -3.32%  "
